$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy the formatting pattern of row 59 down across rows 60-80 ---
$ws.Range("A59:H59").Copy() | Out-Null
$ws.Range("A60:H80").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 2: fix up the "Probleme rencontre" (D) column style: most of the new rows
#     use the wrapped date/time-formatted style (like D2), only D61:D64 keep the
#     plain wrapped style inherited from row 59 ---
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D60").PasteSpecial(-4122) | Out-Null
$ws.Range("D65:D80").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 3: fill in the journal entries for rows 60-80 ---
# Row 60
$ws.Range("A60").Value = 45435
$ws.Range("B60").Value = 'Implementation'
$ws.Range("C60").Value = 'Ajout d''une classe de préférence afin de stocker les tags'
$ws.Range("E60").Value = 0.36805555555555558
$ws.Range("F60").Value = 0.39930555555555558
$ws.Range("G60").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'
$ws.Rows.Item(60).RowHeight = 30

# Row 61
$ws.Range("A61").Value = 45435
$ws.Range("B61").Value = 'Implementation'
$ws.Range("C61").Value = 'Affichage des tags'
$ws.Range("E61").Value = 0.40972222222222227
$ws.Range("F61").Value = 0.50347222222222221
$ws.Range("G61").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 62
$ws.Range("A62").Value = 45435
$ws.Range("B62").Value = 'Implementation'
$ws.Range("C62").Value = 'Création d''objectif'
$ws.Range("E62").Value = 0.5625
$ws.Range("F62").Value = 0.60416666666666663
$ws.Range("G62").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 63
$ws.Range("A63").Value = 45435
$ws.Range("B63").Value = 'Implementation'
$ws.Range("C63").Value = 'Affichage des objectif'
$ws.Range("E63").Value = 0.60416666666666663
$ws.Range("F63").Value = 0.62847222222222221
$ws.Range("G63").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 64
$ws.Range("A64").Value = 45435
$ws.Range("B64").Value = 'Documentation'
$ws.Range("C64").Value = 'Classes mises a jours'
$ws.Range("E64").Value = 0.63888888888888895
$ws.Range("F64").Value = 0.70486111111111116
$ws.Range("G64").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 65
$ws.Range("A65").Value = 45436
$ws.Range("B65").Value = 'Documentation'
$ws.Range("C65").Value = 'Résolution de bug avec mes variables TimeMillis'
$ws.Range("E65").Value = 0.36805555555555558
$ws.Range("F65").Value = 0.39930555555555558
$ws.Range("G65").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 66
$ws.Range("A66").Value = 45436
$ws.Range("B66").Value = 'Implementation'
$ws.Range("C66").Value = 'Ajout de statut'
$ws.Range("E66").Value = 0.40972222222222227
$ws.Range("F66").Value = 0.50347222222222221
$ws.Range("G66").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 67
$ws.Range("A67").Value = 45436
$ws.Range("B67").Value = 'Implementation'
$ws.Range("C67").Value = 'Changement de couleur de police'
$ws.Range("E67").Value = 0.5625
$ws.Range("F67").Value = 0.60416666666666663
$ws.Range("G67").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 68
$ws.Range("A68").Value = 45436
$ws.Range("B68").Value = 'Implementation'
$ws.Range("C68").Value = 'La base de donnée a été mise a jour (version 3)'
$ws.Range("E68").Value = 0.60416666666666663
$ws.Range("F68").Value = 0.62847222222222221
$ws.Range("G68").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 69
$ws.Range("A69").Value = 45436
$ws.Range("B69").Value = 'Documentation'
$ws.Range("C69").Value = 'Rédaction  de  commentaire, correction de bug mineurs'
$ws.Range("E69").Value = 0.63888888888888895
$ws.Range("F69").Value = 0.70486111111111116
$ws.Range("G69").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'
$ws.Rows.Item(69).RowHeight = 30

# Row 70
$ws.Range("A70").Value = 45439
$ws.Range("B70").Value = 'Implementation'
$ws.Range("C70").Value = 'Ajout de la fonction d''annuler un objectif'
$ws.Range("E70").Value = 0.5625
$ws.Range("F70").Value = 0.62847222222222221
$ws.Range("G70").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 71
$ws.Range("A71").Value = 45440
$ws.Range("B71").Value = 'Implementation'
$ws.Range("C71").Value = 'Dernier commit d''implementation, correction de bug mineur'
$ws.Range("E71").Value = 0.33333333333333331
$ws.Range("F71").Value = 0.39930555555555558
$ws.Range("G71").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'
$ws.Rows.Item(71).RowHeight = 30

# Row 72
$ws.Range("A72").Value = 45440
$ws.Range("B72").Value = 'Implementation'
$ws.Range("C72").Value = 'Documentation de la création de la base donnée'
$ws.Range("E72").Value = 0.40972222222222227
$ws.Range("F72").Value = 0.50347222222222221
$ws.Range("G72").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 73
$ws.Range("A73").Value = 45440
$ws.Range("B73").Value = 'Documentation'
$ws.Range("C73").Value = 'Documentation de l''implementation'
$ws.Range("D73").Value = 'Ajout d’un objectif dans la base de données , Affichage des objectifs, Ajout de statut'
$ws.Range("E73").Value = 0.5625
$ws.Range("F73").Value = 0.60416666666666663
$ws.Range("G73").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'
$ws.Rows.Item(73).RowHeight = 30

# Row 74
$ws.Range("A74").Value = 45440
$ws.Range("B74").Value = 'Documentation'
$ws.Range("C74").Value = 'Documentation de l''implementation'
$ws.Range("D74").Value = 'Modification d’un objectif, Tri de l’affichage des objectifs, Annulation d’une notification'
$ws.Range("E74").Value = 0.60416666666666663
$ws.Range("F74").Value = 0.62847222222222221
$ws.Range("G74").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'
$ws.Rows.Item(74).RowHeight = 30

# Row 75
$ws.Range("A75").Value = 45440
$ws.Range("B75").Value = 'Documentation'
$ws.Range("C75").Value = 'Journal de travail et pdf mis a jour'
$ws.Range("D75").Value = 'uploadé sur github'
$ws.Range("E75").Value = 0.63888888888888895
$ws.Range("F75").Value = 0.70486111111111116
$ws.Range("G75").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 76
$ws.Range("A76").Value = 45442
$ws.Range("B76").Value = 'Documentation'
$ws.Range("C76").Value = 'Finalisation du dossier de projet'
$ws.Range("E76").Value = 0.36805555555555558
$ws.Range("F76").Value = 0.39930555555555558
$ws.Range("G76").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 77
$ws.Range("A77").Value = 45442
$ws.Range("B77").Value = 'Documentation'
$ws.Range("C77").Value = 'Finalisation du dossier de projet'
$ws.Range("E77").Value = 0.40972222222222227
$ws.Range("F77").Value = 0.45833333333333331
$ws.Range("G77").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 78
$ws.Range("A78").Value = 45442
$ws.Range("B78").Value = 'Documentation'
$ws.Range("C78").Value = 'Push d''un patch qui corrige trois bug mineur'
$ws.Range("D78").Value = 'découvert lors de la phase de test'
$ws.Range("E78").Value = 0.45833333333333331
$ws.Range("F78").Value = 0.47916666666666669
$ws.Range("G78").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 79
$ws.Range("A79").Value = 45442
$ws.Range("B79").Value = 'Documentation'
$ws.Range("C79").Value = 'Finalisation du dossier de projet'
$ws.Range("E79").Value = 0.47916666666666669
$ws.Range("F79").Value = 0.51041666666666663
$ws.Range("G79").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# Row 80
$ws.Range("A80").Value = 45442
$ws.Range("B80").Value = 'Documentation'
$ws.Range("C80").Value = 'Préparation et rendu du TPI'
$ws.Range("E80").Value = 0.5625
$ws.Range("F80").Value = 0.60069444444444442
$ws.Range("G80").Formula = '=SUM(Tableau1[[#This Row],[Heure de fin]]-Tableau1[[#This Row],[heure de début2]])'

# --- Step 4: grow the Excel table (Tableau1) so it covers the new rows ---
$lo = $ws.ListObjects.Item("Tableau1")
$lo.Resize($ws.Range("A1:H80")) | Out-Null

# --- Step 5: extend the trailing blank C/D rows (was 81-118, now 81-125) ---
$ws.Range("C118:D118").Copy() | Out-Null
$ws.Range("C119:D125").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Step 6: selection + page setup scale ---
$ws.Range("D82").Select() | Out-Null
$ws.PageSetup.Zoom = 60

# --- Step 7: recalc so cached formula values (I2 total, G column, etc.) are fresh ---
$excel.CalculateFullRebuild() | Out-Null
